$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Section 1 test results (K18:K33) ---
$ws.Range("K18").Value = 0.7
$ws.Range("K19").Value = 0.7
$ws.Range("K20").Value = 0.7
$ws.Range("K21").Value = 0.7
$ws.Range("K22").Value = 0.7
$ws.Range("K23").Value = 0.7
$ws.Range("K24").Value = 0.7
$ws.Range("K25").Value = 0.8
$ws.Range("K26").Value = 0.8
$ws.Range("K27").Value = 0.7
$ws.Range("K28").Value = 0.7
$ws.Range("K29").Value = 0.6
$ws.Range("K30").Value = 0.7
$ws.Range("K31").Value = 0.8
$ws.Range("K32").Value = 0.8
$ws.Range("K33").Value = 0.8

# --- Section 1 continued test results (K40:K47) ---
$ws.Range("K40").Value = 0.3
$ws.Range("K41").Value = 0.3
$ws.Range("K42").Value = 0.3
$ws.Range("K43").Value = 0.3
$ws.Range("K44").Value = 0.3
$ws.Range("K45").Value = 0.3
$ws.Range("F46").Value = 30
$ws.Range("K46").Value = 0.3
$ws.Range("K47").Value = 0.3

# --- Final decision (M57:M59) ---
$ws.Range("M57").Value = "y"
$ws.Range("M58").Value = "y"
$ws.Range("M59").Value = "y"

# --- Section 2 result ---
$ws.Range("B60").Value = "PASS"

# --- Section 4 (other tests) measurements ---
$ws.Range("K66").Value = 98.8
$ws.Range("K67").Value = 9.91
$ws.Range("C70").Value = 0.018
$ws.Range("C71").Value = 0.019

# --- Tester & test date ---
$ws.Range("B80").Value = "Brian / Frank"
$ws.Range("F80").Value = "9/26/2014"
